$d = $word.ActiveDocument

# --- Text for the new "Related Work" body paragraph --------------------------------
$text1 = "The closes related system to GRAPH/Z is Pregel, which it was inspired by. In most of our work, however, we compare GRAPH/Z to Graphlab, which is another high performance graph processing framework. Graphlab uses a similar paradigm to GRAPH/Z, but allows a vertex to access data that is not in a message to the vertex."
$text2 = "Another less similar but still relevant work is Hadoop, which follows the MapReduce paradigm. Hadoop is frequently used to process large graphs. In fact, GRAPH/Z and Pregel computations can be expressed as a series of chained MapReduce functions. Hadoop has largely been replaced by Apache Spark, which is faster in some cases."

# --- Locate the "Related Work" heading paragraph and the paragraph just before it ----
$headingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.StartsWith("Related Work")) {
        $headingIndex = $i
    }
}

$beforeHeadingIndex = $headingIndex - 1
$beforePara = $d.Paragraphs.Item($beforeHeadingIndex)

# 1. Insert a new paragraph right after the paragraph preceding the heading (i.e. the
#    "Problem" body text). Because that paragraph has no explicit style override, the
#    freshly inserted paragraph also comes out with the default (Normal) formatting and
#    no <w:pPr> element, matching a plain body paragraph.
$beforePara.Range.InsertParagraphAfter()
$newParaIndex = $beforeHeadingIndex + 1
$newPara = $d.Paragraphs.Item($newParaIndex)

# 2. Populate it with the two sentences, separated by a manual line break, using
#    distinct runs (text, break, text) just like the rest of the document does.
$insPos = $newPara.Range.Start
$r1 = $d.Range($insPos, $insPos)
$r1.InsertAfter($text1)

$insPos2 = $insPos + $text1.Length
$r2 = $d.Range($insPos2, $insPos2)
$r2.InsertBreak(6)

$insPos3 = $insPos2 + 1
$r3 = $d.Range($insPos3, $insPos3)
$r3.InsertAfter($text2)

# 3. The "Related Work" heading paragraph is now immediately after our new paragraph;
#    cut it and paste it back in right before the new paragraph to restore the correct
#    order (heading, then the new body text).
$headingParaIndex = $newParaIndex + 1
$headingPara = $d.Paragraphs.Item($headingParaIndex)
$headingPara.Range.Cut()

$insertPos = $d.Paragraphs.Item($newParaIndex).Range.Start
$target = $d.Range($insertPos, $insertPos)
$target.Paste()

# 4. Cutting the heading paragraph drops the bookmark that marked it; recreate it at
#    the start of the (moved) heading paragraph.
$restoredHeadingPara = $d.Paragraphs.Item($newParaIndex)
$bmRng = $d.Range($restoredHeadingPara.Range.Start, $restoredHeadingPara.Range.Start)
$d.Bookmarks.Add("related-work", $bmRng)

Write-Output "Related Work section added"
